$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph - this marks the start of the
# block to remove (together with the blank paragraph immediately before it).
$startRange = $d.Content
$startRange.Find.ClearFormatting()
$startFound = $startRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the trailing "(c) 2020 ... Creative Commons Attribution" paragraph -
# this marks the end of the block to remove.
$endRange = $d.Content
$endRange.Find.ClearFormatting()
$endFound = $endRange.Find.Execute( `
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($startFound -and $endFound) {
    # Extend one character to the left to also swallow the empty paragraph
    # that precedes "Ver no Jupiter ..." (its trailing paragraph mark), and
    # one character past the end of the copyright line to swallow its own
    # paragraph mark too - this removes all three paragraphs in one shot.
    $deleteStart = $startRange.Start - 1
    $deleteEnd = $endRange.End + 1

    $toDelete = $d.Range($deleteStart, $deleteEnd)
    $toDelete.Delete()
}
